# CryCompanywiseStockReport_1.xlsx data-correction pass.
# Stock quantities (col F) were re-reconciled for a number of line items;
# the corresponding stock value (col G = Rate(D) * Qty(F)) and every
# affected "Sub Total:" / "Grand Total:" row (col B) are updated to match.
# A few adjacent line-item rows were also re-sequenced (their B/C/E/F/G
# data swapped while the Sr.No. in col A stays put).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F19").Value = 2
$ws.Range("G19").Value = 61.48

$ws.Range("F23").Value = 12
$ws.Range("G23").Value = 491.88

$ws.Range("F27").Value = 58
$ws.Range("G27").Value = 2079.88

$ws.Range("F29").Value = 73
$ws.Range("G29").Value = 3739.79

$ws.Range("F31").Value = 50
$ws.Range("G31").Value = 1332

$ws.Range("B34").Value = 59930.36

$ws.Range("F41").Value = 212
$ws.Range("G41").Value = 40892.68

$ws.Range("F44").Value = 33
$ws.Range("G44").Value = 1165.56

$ws.Range("F46").Value = 58
$ws.Range("G46").Value = 2156.44

$ws.Range("F52").Value = 52
$ws.Range("G52").Value = 3068

$ws.Range("F55").Value = 127
$ws.Range("G55").Value = 7081.52

$ws.Range("F57").Value = 50
$ws.Range("G57").Value = 1766

$ws.Range("F61").Value = 233
$ws.Range("G61").Value = 60750.09

$ws.Range("B66").Value = 208502.58

$ws.Range("F106").Value = 182
$ws.Range("G106").Value = 20416.76

$ws.Range("F114").Value = 62
$ws.Range("G114").Value = 2897.88

$ws.Range("B123").Value = 74823.99000000001

$ws.Range("B126").Value = 65258
$ws.Range("F126").Value = 0
$ws.Range("G126").Value = 0

$ws.Range("B127").Value = 64196
$ws.Range("F127").Value = 1
$ws.Range("G127").Value = 32143.58

$ws.Range("F146").Value = 18
$ws.Range("G146").Value = 1182.06

$ws.Range("B147").Value = 22046.32

$ws.Range("F152").Value = 134
$ws.Range("G152").Value = 6393.14

$ws.Range("F154").Value = 303
$ws.Range("G154").Value = 10102.02

$ws.Range("B155").Value = 37824.01

$ws.Range("F182").Value = 20
$ws.Range("G182").Value = 1790.8

$ws.Range("F186").Value = 18
$ws.Range("G186").Value = 779.04

$ws.Range("B193").Value = 66162.39999999999

$ws.Range("F206").Value = 72
$ws.Range("G206").Value = 4665.6

$ws.Range("B208").Value = 4712.09

$ws.Range("F210").Value = 130
$ws.Range("G210").Value = 7073.3

$ws.Range("F211").Value = 8
$ws.Range("G211").Value = 816.08

$ws.Range("F212").Value = 66
$ws.Range("G212").Value = 5879.94

$ws.Range("F217").Value = 43
$ws.Range("G217").Value = 3194.9

$ws.Range("B218").Value = 81027.64999999999

$ws.Range("F222").Value = 915
$ws.Range("G222").Value = 16927.5

$ws.Range("B229").Value = 29054.42

$ws.Range("F252").Value = 75
$ws.Range("G252").Value = 6686.25

$ws.Range("F267").Value = 138
$ws.Range("G267").Value = 5862.24

$ws.Range("F288").Value = 5
$ws.Range("G288").Value = 2775.15

$ws.Range("B290").Value = 66194
$ws.Range("C290").Value = "HIM-Total Care Baby Pants Diapers-M-9s"
$ws.Range("F290").Value = 27
$ws.Range("G290").Value = 2313.36

$ws.Range("B291").Value = 64983
$ws.Range("C291").Value = "HIM-TOTAL CARE BABY PANTS DIAPERS-M-9S"
$ws.Range("F291").Value = 6
$ws.Range("G291").Value = 514.08

$ws.Range("B295").Value = 122746.89

$ws.Range("B304").Value = 63520
$ws.Range("E304").Value = 153.4
$ws.Range("F304").Value = 39
$ws.Range("G304").Value = 5626.92

$ws.Range("B305").Value = 55373
$ws.Range("E305").Value = 163.62
$ws.Range("F305").Value = -94
$ws.Range("G305").Value = -13562.32

$ws.Range("F324").Value = 41
$ws.Range("G324").Value = 7024.53

$ws.Range("B328").Value = -5761.98

$ws.Range("F349").Value = 150
$ws.Range("G349").Value = 11194.5

$ws.Range("B356").Value = 79215.25

$ws.Range("F358").Value = 45
$ws.Range("G358").Value = 10359.9

$ws.Range("B363").Value = 77931.75999999999

$ws.Range("F368").Value = 59
$ws.Range("G368").Value = 1889.18

$ws.Range("F370").Value = 229
$ws.Range("G370").Value = 38011.71

$ws.Range("F371").Value = 65
$ws.Range("G371").Value = 9767.549999999999

$ws.Range("B372").Value = 63231.39

$ws.Range("F387").Value = 435
$ws.Range("G387").Value = 42021

$ws.Range("B389").Value = 58962.46

$ws.Range("F394").Value = 154
$ws.Range("G394").Value = 7980.28

$ws.Range("F396").Value = 133
$ws.Range("G396").Value = 3388.84

$ws.Range("F402").Value = 50
$ws.Range("G402").Value = 1715.5

$ws.Range("F416").Value = 74
$ws.Range("G416").Value = 2172.64

$ws.Range("B417").Value = 173820

$ws.Range("F455").Value = 45
$ws.Range("G455").Value = 10000.35

$ws.Range("B458").Value = 99854.67

$ws.Range("B479").Value = 53319
$ws.Range("E479").Value = 310.64
$ws.Range("F479").Value = -6
$ws.Range("G479").Value = -1643.52

$ws.Range("B480").Value = 64810
$ws.Range("E480").Value = 291.22
$ws.Range("F480").Value = 0
$ws.Range("G480").Value = 0

$ws.Range("B496").Value = 64833
$ws.Range("E496").Value = 34.9
$ws.Range("F496").Value = 88
$ws.Range("G496").Value = 2889.04

$ws.Range("B497").Value = 60025
$ws.Range("E497").Value = 37.22
$ws.Range("F497").Value = -98
$ws.Range("G497").Value = -3217.34

$ws.Range("B506").Value = 60022
$ws.Range("E506").Value = 37.22
$ws.Range("F506").Value = -113
$ws.Range("G506").Value = -3709.79

$ws.Range("B507").Value = 64830
$ws.Range("E507").Value = 34.9
$ws.Range("F507").Value = 84
$ws.Range("G507").Value = 2757.72

$ws.Range("F511").Value = 249
$ws.Range("G511").Value = 24867.63

$ws.Range("F513").Value = 222
$ws.Range("G513").Value = 10982.34

$ws.Range("F517").Value = 149
$ws.Range("G517").Value = 8831.23

$ws.Range("F519").Value = 416
$ws.Range("G519").Value = 22830.08

$ws.Range("F524").Value = 22
$ws.Range("G524").Value = 1943.48

$ws.Range("B525").Value = 129159.96

$ws.Range("F529").Value = 122
$ws.Range("G529").Value = 4039.42

$ws.Range("F530").Value = 20
$ws.Range("G530").Value = 863.6

$ws.Range("F534").Value = 130
$ws.Range("G534").Value = 5688.8

$ws.Range("B535").Value = 24634.73

$ws.Range("F543").Value = 4
$ws.Range("G543").Value = 1481.08

$ws.Range("F544").Value = 42
$ws.Range("G544").Value = 2599.8

$ws.Range("F555").Value = 5
$ws.Range("G555").Value = 1270.05

$ws.Range("B556").Value = 50316.59

$ws.Range("F558").Value = 204
$ws.Range("G558").Value = 24857.4

$ws.Range("B561").Value = 29347.6

$ws.Range("F566").Value = 4
$ws.Range("G566").Value = 1305.84

$ws.Range("B573").Value = 26975.99

$ws.Range("F601").Value = 11
$ws.Range("G601").Value = 1424.83

$ws.Range("B603").Value = 6519.88

$ws.Range("F605").Value = 186
$ws.Range("G605").Value = 24756.6

$ws.Range("B607").Value = 25161.63

$ws.Range("F620").Value = 366
$ws.Range("G620").Value = 28763.94

$ws.Range("F622").Value = 487
$ws.Range("G622").Value = 50117.17

$ws.Range("B628").Value = 212410.11

$ws.Range("F696").Value = 5
$ws.Range("G696").Value = 2183.5

$ws.Range("B713").Value = 69496.44

$ws.Range("B718").Value = 2792644.06

$ws.Range("B719").Value = 2792644.06
